# Workbook: DaySale_2025-06-03_00-00.xlsx
# Change: a new product line ("مبرد قدم") was added to the shortage report,
# pushing the previous last line ("محلول ملح") down one row, bumping the
# running total, and refreshing the generation timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row right before the old totals row (row 34), shifting the
#    totals row -> 35 and the footer row -> 36.
$ws.Rows(34).Insert()

# 2) Clone row 33's cell formatting (and, for now, its values) down into the
#    freshly inserted row 34 so the new product row has the right look.
$ws.Range("A33:Q33").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3) Re-create the merges for the new row 34 (mirrors row 33's merge layout).
$ws.Range("A34:B34").Merge()
$ws.Range("C34:G34").Merge()
$ws.Range("H34:K34").Merge()
$ws.Range("L34:M34").Merge()
$ws.Range("N34:O34").Merge()

# 4) Row 34 now holds a duplicate of the old row 33 ("محلول ملح") - just fix
#    up its running index (27 -> 28); the rest of the copied data is already
#    correct for that product.
$ws.Range("A34").Value = 28

# 5) Row 33 becomes the new product line ("مبرد قدم"), reusing the same
#    numbers ("3:0" / "1" / "20.00" / "20.0000") seen elsewhere in the sheet.
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "'مبرد قدم"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "'3:0"
$ws.Range("L33").NumberFormat = "@"
$ws.Range("L33").Value = "'1"
$ws.Range("N33").NumberFormat = "@"
$ws.Range("N33").Value = "'20.00"
$ws.Range("P33").NumberFormat = "@"
$ws.Range("P33").Value = "'20.0000"

# 6) Row heights: row 33 keeps its height, the new product row 34 takes the
#    height the totals row used to have, the totals row (now 35) gets the
#    standard product-row height, and the footer row (now 36) is untouched.
$ws.Rows(33).RowHeight = 24.75
$ws.Rows(34).RowHeight = 25.5
$ws.Rows(35).RowHeight = 24.75
$ws.Rows(36).RowHeight = 16.5

# 7) Bump the grand total shown under column P (was 1229.15, +20.00 for the
#    new line == 1249.15) - this lives in the now-shifted totals row 35.
$ws.Range("P35").Value = 1249.1500000000001

# 8) Refresh the "generated at" timestamp in the footer (now row 36).
$ws.Range("A36").Value = "Tuesday, 3 June, 2025 11:51 AM"
